# "agregue un link al repositorio"
#
# 1. Extend the paragraph that ends with "...capturas de pantalla." with a
#    new sentence that links to the GitHub repository (pushing the existing
#    paragraph break after the new text).
# 2. Append a trailing (unformatted) space run to the "Comenzamos con Git
#    status:" paragraph.
# 3. Remove the (hidden) "_GoBack" bookmark that used to sit in the middle
#    of "...dejamos una captura de pantalla de ejemplo."
# 4. Drop the stale <w:lastRenderedPageBreak/> that preceded "El Git Push...".

$d = $word.ActiveDocument

# --- 1. Insert the repository link sentence -------------------------------
# Borrow the sz=32/szCs=32 run formatting of the text that immediately
# precedes the insertion point, via FormattedText, so the new sentence gets
# the same explicit run properties as the rest of the paragraph.
$tmpl = $d.Content
$tmpl.Find.Execute("capturas de pantalla.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tmplLen = $tmpl.End - $tmpl.Start
$ft = $tmpl.FormattedText

$dest = $d.Content
$dest.Find.Execute("y los diversos comandos demostrados con las capturas de pantalla.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dest.Collapse(0)
$insertStart = $dest.Start
$dest.FormattedText = $ft
$newRng = $d.Range($insertStart, $insertStart + $tmplLen)
$newRng.Text = " A continuación, les dejamos el link hacia nuestro repositorio: https://github.com/tomasloray/Trabajo-Especial-git-2018---reentrega"

# --- 2. Trailing space after "Comenzamos con Git status:" -----------------
$rng2 = $d.Content
$rng2.Find.Execute("Comenzamos con Git status:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)
$rng2.InsertAfter(" ")

# --- 3. Remove the stray "_GoBack" bookmark --------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 4. Remove the stale lastRenderedPageBreak before "El Git Push" -------
$pb = $d.Content
$pb.Find.Execute("El Git Push", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pbStart = $pb.Start
$pbNarrow = $d.Range($pbStart, $pbStart + 3)
$pbNarrow.Text = "PLACEHOLDER_EL_TAG"
$pbFix = $d.Content
$pbFix.Find.Execute("PLACEHOLDER_EL_TAG", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pbFix.Text = "El "

Write-Host "done"
